# Update "想去人数" (column F) values on the "展览" and "全部类型" worksheets
# to reflect newly generated output (commit: Update gh-pages to output
# generated at 456a3b4).

$wb = $excel.ActiveWorkbook

# Row -> new F value for sheet "展览"
$exhibitionUpdates = @{
    2  = 610
    3  = 495
    6  = 14262
    7  = 16216
    17 = 36
    18 = 94
    24 = 6508
    27 = 1111
    29 = 5690
    31 = 143
    32 = 169
    33 = 4718
    34 = 15
}

# Row -> new F value for sheet "全部类型"
$allTypesUpdates = @{
    2  = 610
    3  = 495
    6  = 14262
    7  = 16216
    17 = 36
    18 = 94
    25 = 6508
    28 = 1111
    31 = 5690
    33 = 143
    34 = 169
    35 = 4718
    36 = 15
}

$wsExhibition = $wb.Worksheets.Item("展览")
foreach ($row in $exhibitionUpdates.Keys) {
    $wsExhibition.Range("F$row").Value = $exhibitionUpdates[$row]
}

$wsAllTypes = $wb.Worksheets.Item("全部类型")
foreach ($row in $allTypesUpdates.Keys) {
    $wsAllTypes.Range("F$row").Value = $allTypesUpdates[$row]
}
